$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '22.450.91'
$ws.Range('E2').Value = '  +0.15%  '
Set-TextValue $ws.Range('D3') '1.573.30'
$ws.Range('E3').Value = '  +0.58%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('E5').Value = '  +0.03%  '
Set-TextValue $ws.Range('D6') '287.97'
$ws.Range('E6').Value = '  +0.69%  '
Set-TextValue $ws.Range('D7') '0.3703'
$ws.Range('E7').Value = '  +1.63%  '
Set-TextValue $ws.Range('D8') '47.22'
$ws.Range('E8').Value = '  -2.59%  '
Set-TextValue $ws.Range('D9') '0.3321'
$ws.Range('E9').Value = '  -0.49%  '
$ws.Range('E10').Value = '  +2.14%  '
Set-TextValue $ws.Range('D11') '0.07500'
$ws.Range('E11').Value = '  +1.08%  '
$ws.Range('E12').Value = '  +0.10%  '
Set-TextValue $ws.Range('D13') '20.78'
$ws.Range('E13').Value = '  -0.29%  '
Set-TextValue $ws.Range('D14') '5.937'
$ws.Range('E14').Value = '  -0.03%  '
Set-TextValue $ws.Range('D15') '6.924'
$ws.Range('E15').Value = '  +0.42%  '
Set-TextValue $ws.Range('D16') '1.563.65'
$ws.Range('E16').Value = '  -0.06%  '
$ws.Range('E17').Value = '  +0.90%  '
Set-TextValue $ws.Range('D18') '88.43'
Set-TextValue $ws.Range('D19') '0.06721'
$ws.Range('E19').Value = '  +0.52%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range('D20') '6.390'
$ws.Range('E20').Value = '  +0.37%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range('D21') '0.9995'
$ws.Range('E21').Value = '  -0.07%  '
Set-TextValue $ws.Range('D22') '16.50'
$ws.Range('E22').Value = '  +2.27%  '
$ws.Range('E23').Value = '  +0.00%  '
Set-TextValue $ws.Range('D24') '22.448.62'
$ws.Range('E24').Value = '  +0.20%  '
Set-TextValue $ws.Range('D25') '2.384'
$ws.Range('E25').Value = '  -1.39%  '
Set-TextValue $ws.Range('D26') '2.639'
$ws.Range('E26').Value = '  +2.90%  '
Set-TextValue $ws.Range('D27') '150.74'
$ws.Range('E27').Value = '  +0.54%  '
$ws.Range('E28').Value = '  +0.92%  '
Set-TextValue $ws.Range('D29') '4.967'
$ws.Range('E29').Value = '  -0.68%  '
Set-TextValue $ws.Range('D30') '125.01'
$ws.Range('E30').Value = '  +1.55%  '
Set-TextValue $ws.Range('D31') '1.741.52'
$ws.Range('E31').Value = '  +0.19%  '
Set-TextValue $ws.Range('D32') '1.095'
$ws.Range('E32').Value = '  +2.70%  '
Set-TextValue $ws.Range('D33') '6.083'
$ws.Range('E33').Value = '  -1.10%  '
Set-TextValue $ws.Range('D34') '1.986'
$ws.Range('E34').Value = '  -0.33%  '
Set-TextValue $ws.Range('D35') '9.895'
$ws.Range('E35').Value = '  +1.89%  '
Set-TextValue $ws.Range('D36') '0.08347'
$ws.Range('E36').Value = '  +1.19%  '
Set-TextValue $ws.Range('D37') '0.02447'
$ws.Range('E37').Value = '  +2.14%  '
Set-TextValue $ws.Range('D38') '1.311'
$ws.Range('E38').Value = '  +0.15%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D39') '0.06391'
$ws.Range('E39').Value = '  +0.11%  '
$ws.Range('B40').Value = 'Algorand'
$ws.Range('C40').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range('D40') '0.2219'
$ws.Range('E40').Value = '  +0.40%  '
Set-TextValue $ws.Range('D41') '5.338'
$ws.Range('E41').Value = '  -0.07%  '
Set-TextValue $ws.Range('D42') '11.40'
$ws.Range('E42').Value = '  +2.08%  '
Set-TextValue $ws.Range('D43') '0.6243'
$ws.Range('E43').Value = '  +2.45%  '
Set-TextValue $ws.Range('D44') '14.03'
$ws.Range('E44').Value = '  +1.72%  '
Set-TextValue $ws.Range('D45') '0.6066'
$ws.Range('E45').Value = '  +5.29%  '
Set-TextValue $ws.Range('D46') '3.775'
$ws.Range('E46').Value = '  +0.39%  '
Set-TextValue $ws.Range('D47') '2.046'
$ws.Range('E47').Value = '  +1.46%  '
Set-TextValue $ws.Range('D48') '124.85'
$ws.Range('E48').Value = '  -0.20%  '
Set-TextValue $ws.Range('D49') '1.208'
$ws.Range('E49').Value = '  -0.65%  '
Set-TextValue $ws.Range('D50') '0.07197'
$ws.Range('E50').Value = '  -0.24%  '
Set-TextValue $ws.Range('D51') '77.21'
$ws.Range('E51').Value = '  +2.53%  '
